$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "23.493.30"
Set-TextValue "E2" "  +0.83%  "
Set-TextValue "D3" "1.646.81"
Set-TextValue "E3" "  +1.51%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "1.002"
Set-TextValue "E5" "  -0.14%  "
Set-TextValue "D6" "302.41"
Set-TextValue "E6" "  -0.20%  "
Set-TextValue "D7" "0.3839"
Set-TextValue "D8" "0.3592"
Set-TextValue "E8" "  +1.72%  "
Set-TextValue "D9" "51.02"
Set-TextValue "D10" "0.08189"
Set-TextValue "E10" "  +1.37%  "
Set-TextValue "D11" "1.231"
Set-TextValue "E11" "  +1.88%  "
Set-TextValue "D12" "1.002"
Set-TextValue "E12" "  +0.00%  "
Set-TextValue "D13" "22.30"
Set-TextValue "E13" "  +1.39%  "
Set-TextValue "D14" "6.458"
Set-TextValue "E14" "  +1.20%  "
Set-TextValue "D15" "7.461"
Set-TextValue "E15" "  +3.42%  "
Set-TextValue "D16" "0.00001224"
Set-TextValue "E16" "  +0.80%  "
Set-TextValue "D17" "1.645.99"
Set-TextValue "E17" "  +1.65%  "
Set-TextValue "D18" "97.49"
Set-TextValue "E18" "  +3.45%  "
Set-TextValue "D19" "0.07003"
Set-TextValue "E19" "  +1.10%  "
Set-TextValue "D20" "6.781"
Set-TextValue "E20" "  +4.26%  "
Set-TextValue "D21" "17.55"
Set-TextValue "E21" "  +2.13%  "
Set-TextValue "E22" "  -0.13%  "
Set-TextValue "D23" "12.65"
Set-TextValue "E23" "  +3.10%  "
Set-TextValue "D24" "23.491.68"
Set-TextValue "E24" "  +0.86%  "
Set-TextValue "D25" "2.499"
Set-TextValue "E25" "  -2.21%  "
Set-TextValue "D26" "3.022"
Set-TextValue "E26" "  -2.68%  "
Set-TextValue "D27" "21.23"
Set-TextValue "E27" "  +1.94%  "
Set-TextValue "D28" "152.87"
Set-TextValue "E28" "  +1.10%  "
Set-TextValue "D29" "5.230"
Set-TextValue "E29" "  -0.44%  "
Set-TextValue "D30" "133.94"
Set-TextValue "E30" "  +1.51%  "
Set-TextValue "D31" "1.829.09"
Set-TextValue "E31" "  +1.51%  "
Set-TextValue "D32" "7.105"
Set-TextValue "E32" "  +10.26%  "
Set-TextValue "D33" "2.244"
Set-TextValue "E33" "  +6.01%  "
Set-TextValue "D34" "12.11"
Set-TextValue "E34" "  +5.71%  "
Set-TextValue "D35" "1.060"
Set-TextValue "E35" "  -0.44%  "
Set-TextValue "D36" "0.02793"
Set-TextValue "E36" "  +2.85%  "
Set-TextValue "D37" "0.2495"
Set-TextValue "E37" "  +1.17%  "
Set-TextValue "D38" "6.090"
Set-TextValue "E38" "  +4.13%  "
Set-TextValue "D39" "0.08763"
Set-TextValue "D40" "0.06985"
Set-TextValue "E40" "  +1.32%  "
Set-TextValue "D41" "13.10"
Set-TextValue "E41" "  +10.03%  "
Set-TextValue "D42" "0.6980"
Set-TextValue "E42" "  +1.59%  "
Set-TextValue "D43" "1.333"
Set-TextValue "E43" "  +1.83%  "
Set-TextValue "D44" "15.83"
Set-TextValue "E44" "  +4.24%  "
Set-TextValue "D45" "0.6513"
Set-TextValue "E45" "  +3.50%  "
Set-TextValue "D47" "2.298"
Set-TextValue "E47" "  +2.29%  "
Set-TextValue "D48" "3.953"
Set-TextValue "E48" "  +0.00%  "
Set-TextValue "D49" "0.07877"
Set-TextValue "E49" "  -0.01%  "
Set-TextValue "D50" "128.04"
Set-TextValue "E50" "  -0.36%  "
Set-TextValue "E51" "  +1.16%  "
